# Fixed the render issues, now can draw level no problem
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the EARNED (column C) values to match the VALUE (column B) for the
# core-feature rows that are now considered complete.
$ws.Range("C6").Value = 0.03
$ws.Range("C7").Value = 0.1
$ws.Range("C9").Value = 0.1
$ws.Range("C10").Value = 0.1
$ws.Range("C11").Value = 0.05

# Match the styling already used on the other completed rows (C2/C3/C4/C5/C8)
$ws.Range("C6").Font.Color = $ws.Range("C8").Font.Color
$ws.Range("C7").Font.Color = $ws.Range("C8").Font.Color
$ws.Range("C9").Font.Color = $ws.Range("C8").Font.Color
$ws.Range("C10").Font.Color = $ws.Range("C8").Font.Color
$ws.Range("C11").Font.Color = $ws.Range("C8").Font.Color

# Update the current selection to reflect where the user ended up working
$ws.Range("G10").Select()
